# Update the "DATA" worksheet: the HOURS value logged for Zerin's
# "Simulation" entry (row 10) changes from 5 to 4, and the active
# selection moves from F10 to C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

$ws.Range("C10").Value = 4

$ws.Activate()
$ws.Range("C10").Select()
